$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to Text format so numeric-looking values
# (e.g. "0.999", "43.46") are stored as text, matching the source data
# which uses inline strings for all price/volume cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "67.500.65"
$ws.Range("E2").Value = "  +1.26%  "
$ws.Range("D3").Value = "3.948.31"
$ws.Range("E3").Value = "  +4.56%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "479.55"
$ws.Range("E5").Value = "  +9.14%  "
$ws.Range("D6").Value = "148.75"
$ws.Range("E6").Value = "  +4.19%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  +1.06%  "
$ws.Range("D9").Value = "0.733"
$ws.Range("E9").Value = "  +0.16%  "
$ws.Range("D10").Value = "0.167"
$ws.Range("E10").Value = "  +10.62%  "
$ws.Range("D11").Value = "0.0000352"
$ws.Range("E11").Value = "  +13.28%  "
$ws.Range("D12").Value = "43.46"
$ws.Range("E12").Value = "  +1.77%  "
$ws.Range("D13").Value = "4.568.57"
$ws.Range("E13").Value = "  +3.94%  "
$ws.Range("D14").Value = "10.46"
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("D15").Value = "14.99"
$ws.Range("E15").Value = "  +1.38%  "
$ws.Range("D16").Value = "3.946.72"
$ws.Range("E16").Value = "  +4.72%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "20.09"
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("E19").Value = "  +2.15%  "
$ws.Range("D20").Value = "67.661.11"
$ws.Range("E20").Value = "  +1.33%  "
$ws.Range("D21").Value = "434.23"
$ws.Range("E21").Value = "  +5.46%  "
$ws.Range("D22").Value = "3.39"
$ws.Range("E22").Value = "  +3.65%  "
$ws.Range("D23").Value = "14.56"
$ws.Range("E23").Value = "  +0.55%  "
$ws.Range("D24").Value = "87.60"
$ws.Range("E24").Value = "  +2.70%  "
$ws.Range("E25").Value = "  +8.04%  "
$ws.Range("D26").Value = "38.66"
$ws.Range("E26").Value = "  +4.95%  "
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +4.60%  "
$ws.Range("D28").Value = "9.78"
$ws.Range("E28").Value = "  +1.67%  "
$ws.Range("D29").Value = "719.82"
$ws.Range("E29").Value = "  -1.55%  "
$ws.Range("E30").Value = "  -0.81%  "
$ws.Range("D31").Value = "13.47"
$ws.Range("E31").Value = "  -2.60%  "
$ws.Range("D32").Value = "2.82"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "42.33"
$ws.Range("E33").Value = "  -0.66%  "
$ws.Range("D34").Value = "0.0₃0842"
$ws.Range("E34").Value = "  +26.23%  "
$ws.Range("E35").Value = "  +3.58%  "
$ws.Range("D36").Value = "0.152"
$ws.Range("E36").Value = "  -3.29%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "5.39"
$ws.Range("E38").Value = "  -3.15%  "
$ws.Range("D39").Value = "0.0476"
$ws.Range("E39").Value = "  -0.15%  "
$ws.Range("E40").Value = "  +5.76%  "
$ws.Range("D41").Value = "0.144"
$ws.Range("E41").Value = "  +3.38%  "
$ws.Range("B42").Value = "TheGraph"
$ws.Range("C42").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D42").Value = "0.338"
$ws.Range("E42").Value = "  +2.74%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "2.23"
$ws.Range("E43").Value = "  +7.31%  "
$ws.Range("E44").Value = "  -0.16%  "
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "2.82"
$ws.Range("E45").Value = "  +6.17%  "
$ws.Range("B46").Value = "LidoDAOToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D46").Value = "3.48"
$ws.Range("E46").Value = "  +5.87%  "
$ws.Range("E47").Value = "  -7.05%  "
$ws.Range("D48").Value = "3.21"
$ws.Range("E48").Value = "  -4.29%  "
$ws.Range("D49").Value = "147.79"
$ws.Range("E49").Value = "  +3.29%  "
$ws.Range("D50").Value = "2.89"
$ws.Range("E50").Value = "  +2.89%  "
$ws.Range("D51").Value = "25.44"
$ws.Range("E51").Value = "  +3.27%  "
